# This script applies a "new week of data" update to the Acelga price sheet.
# A new record (row 60, date 2021-08-05) is inserted at the top of the
# date-ordered block that starts at row 60. Every existing record in rows
# 60-143 shifts down by one row (row 60 -> 61, 61 -> 62, ... 142 -> 143),
# and the record that used to be on row 143 now lands on the brand new
# row 144 at the bottom of the sheet.
#
# Only the columns that actually vary between records need to move:
#   D (Fecha), J (Volumen), K (Precio minimo), L (Precio maximo),
#   M (Precio promedio ponderado), O (Origen), P (Precio $/Kg)
# All the other columns (A,B,C,E,F,G,H,I,N,Q,R) are constant across every
# record in this sheet, so the new row 144 can simply be populated with
# those same constant values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 60
$lastRow = 143
$newLastRow = 144

# Populate the constant columns for the brand new row at the bottom.
$ws.Cells.Item($newLastRow, 1).Value = 5
$ws.Cells.Item($newLastRow, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item($newLastRow, 3).Value = "Maule"
$ws.Cells.Item($newLastRow, 5).Value = 7
$ws.Cells.Item($newLastRow, 6).Value = 100112009
$ws.Cells.Item($newLastRow, 7).Value = "Acelga"
$ws.Cells.Item($newLastRow, 8).Value = "Sin especificar"
$ws.Cells.Item($newLastRow, 9).Value = "Primera"
$ws.Cells.Item($newLastRow, 14).Value = "`$/docena de atados (4 kilos)"
$ws.Cells.Item($newLastRow, 17).Value = 4
$ws.Cells.Item($newLastRow, 18).Value = "Hortaliza"

# The Fecha (date) column uses a custom date number format throughout the
# sheet; apply the same format to the new row before writing its date so a
# duplicate style entry isn't created.
$ws.Cells.Item($newLastRow, 4).NumberFormat = $ws.Cells.Item($lastRow, 4).NumberFormat()

# Shift the variable columns down by one row, starting from the bottom so
# that the source row for each copy hasn't been overwritten yet.
for ($r = $lastRow; $r -ge $firstRow; $r--) {
    $destRow = $r + 1

    $dVal = $ws.Cells.Item($r, 4).Value()
    $jVal = $ws.Cells.Item($r, 10).Value()
    $kVal = $ws.Cells.Item($r, 11).Value()
    $lVal = $ws.Cells.Item($r, 12).Value()
    $mVal = $ws.Cells.Item($r, 13).Value()
    $oVal = $ws.Cells.Item($r, 15).Value()
    $pVal = $ws.Cells.Item($r, 16).Value()

    $ws.Cells.Item($destRow, 4).Value = $dVal
    $ws.Cells.Item($destRow, 10).Value = $jVal
    $ws.Cells.Item($destRow, 11).Value = $kVal
    $ws.Cells.Item($destRow, 12).Value = $lVal
    $ws.Cells.Item($destRow, 13).Value = $mVal
    $ws.Cells.Item($destRow, 15).Value = $oVal
    $ws.Cells.Item($destRow, 16).Value = $pVal
}

# Finally, write the brand new record into row 60.
$ws.Cells.Item($firstRow, 4).Value = 44413
$ws.Cells.Item($firstRow, 10).Value = 400
$ws.Cells.Item($firstRow, 11).Value = 2000
$ws.Cells.Item($firstRow, 12).Value = 2000
$ws.Cells.Item($firstRow, 13).Value = 2000
$ws.Cells.Item($firstRow, 15).Value = "Región del Maule"
$ws.Cells.Item($firstRow, 16).Value = 500
